$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 105
$ws.Range("F4").Value = 8257
$ws.Range("F6").Value = 109
$ws.Range("F7").Value = 7210
$ws.Range("F8").Value = 1140
$ws.Range("F9").Value = 564
$ws.Range("F10").Value = 497
$ws.Range("F11").Value = 719
$ws.Range("F13").Value = 164
$ws.Range("F17").Value = 105
$ws.Range("F18").Value = 11844
$ws.Range("F19").Value = 103
$ws.Range("F21").Value = 140
$ws.Range("F22").Value = 2353
$ws.Range("F24").Value = 3322
$ws.Range("F27").Value = 2798
$ws.Range("F28").Value = 108
$ws.Range("F29").Value = 29
$ws.Range("F31").Value = 3295
$ws.Range("F33").Value = 2409
$ws.Range("F35").Value = 1661
$ws.Range("F38").Value = 5897
$ws.Range("F39").Value = 1252
$ws.Range("F40").Value = 15
$ws.Range("F42").Value = 192
$ws.Range("F43").Value = 1121
$ws.Range("F44").Value = 1106
$ws.Range("F45").Value = 1084
$ws.Range("F46").Value = 1552
$ws.Range("F47").Value = 10
$ws.Range("F48").Value = 104
$ws.Range("F49").Value = 1134

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 917
$ws.Range("F20").Value = 70

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 278
$ws.Range("F3").Value = 421
$ws.Range("E4").Value = '2024.07.27 10:00-07.28 22:00'
$ws.Range("F4").Value = 7

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 105
$ws.Range("F4").Value = 278
$ws.Range("F5").Value = 421
$ws.Range("F8").Value = 8257
$ws.Range("F10").Value = 109
$ws.Range("F11").Value = 7210
$ws.Range("F12").Value = 7210
$ws.Range("F13").Value = 1140
$ws.Range("F14").Value = 564
$ws.Range("F15").Value = 497
$ws.Range("F16").Value = 719
$ws.Range("F18").Value = 164
$ws.Range("F20").Value = 105
$ws.Range("F22").Value = 11844
$ws.Range("F23").Value = 103
$ws.Range("F25").Value = 140
$ws.Range("F26").Value = 2353
$ws.Range("F27").Value = 2353
$ws.Range("F28").Value = 3322
$ws.Range("F29").Value = 2798
$ws.Range("F30").Value = 108
$ws.Range("F31").Value = 29
$ws.Range("F33").Value = 3295
$ws.Range("F36").Value = 2409
$ws.Range("F38").Value = 1661
$ws.Range("F40").Value = 5897
$ws.Range("F41").Value = 70
$ws.Range("F43").Value = 1252
$ws.Range("F45").Value = 192
$ws.Range("F46").Value = 1121
$ws.Range("F47").Value = 1106
$ws.Range("F48").Value = 1084
$ws.Range("F49").Value = 1552
$ws.Range("F50").Value = 104
$ws.Range("F51").Value = 1134

